$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: Correspond Handoff Datetime (E) and Correspond Handback DateTime (H)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2:E3").Value = "2016-03-24 03:12:32"
$wsZhCn.Range("H2:H3").Value = "2016-03-24 03:13:14"

# "de-de" sheet: Correspond Handoff Datetime (E) and Correspond Handback DateTime (H)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2:E3").Value = "2016-03-24 03:12:40"
$wsDeDe.Range("H2:H3").Value = "2016-03-24 03:13:14"
